$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that currently sits right
#    before the "Things we considered" heading (it is being moved).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Locate the sentence that needs to be re-split / re-coloured.
# ------------------------------------------------------------------
$sentence = "It was clear that we will use convolution layers to capture spatial patterns and reduce the complexity of the model."
$find = $d.Content
$found = $find.Find.Execute($sentence, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$sentenceStart = $find.Start

$part1 = "It was clear that we will use "
$part2 = "convolution layers "
$part3 = "to "
$part4 = "capture "
$part5 = "spatial patterns and reduce the complexity of the model."

$p1End = $sentenceStart + $part1.Length
$p2End = $p1End + $part2.Length
$p3End = $p2End + $part3.Length
$p4End = $p3End + $part4.Length
$p5End = $p4End + $part5.Length

# ------------------------------------------------------------------
# 3) Strip the red colouring from the "to " run by deleting it and
#    retyping plain text in its place (the fresh text carries no
#    explicit run formatting).
# ------------------------------------------------------------------
$toRange = $d.Range($p2End, $p3End)
$toRange.Delete()
$toInsertionPoint = $d.Range($p2End, $p2End)
$toInsertionPoint.InsertAfter($part3)

# ------------------------------------------------------------------
# 4) Force genuine run boundaries at each of the split points so the
#    identically-formatted neighbouring text does not get re-merged
#    into a single run. A temporary, zero-length bookmark dropped at
#    a boundary and then removed leaves the run split behind.
# ------------------------------------------------------------------
function Force-Split([int]$pos) {
    $name = "ZzTempSplit"
    $r = $d.Range($pos, $pos)
    $d.Bookmarks.Add($name, $r)
    $d.Bookmarks($name).Delete()
}

Force-Split($p1End)
Force-Split($p2End)
Force-Split($p3End)
Force-Split($p4End)

# ------------------------------------------------------------------
# 5) Re-add the "_GoBack" bookmark around "convolution layers to
#    capture " (its new location per the edit).
# ------------------------------------------------------------------
$bmRange = $d.Range($p1End, $p4End)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output $d.Range($sentenceStart, $p5End).Text
